# Daily attendance processing - 2026-01-28 22:41:41
# Normalizes the ordering of names/emails in the "Recorded By" column (G)
# for the "Session Analysis Results" sheet. A handful of distinct
# "Recorded By" combinations had their constituent names listed in an
# inconsistent order; this pass rewrites each occurrence to the corrected
# ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact-match lookup of old "Recorded By" text -> corrected text.
$map = @{
    "System, backup@backdoor.com, system" = "system, System, backup@backdoor.com"
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "backup@backdoor.com, System"         = "System, backup@backdoor.com"
    "System, admin@admin.com"             = "admin@admin.com, System"
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com"
}

$lastRow = $ws.UsedRange.Rows.Count
$col = 7   # column G = "Recorded By"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $current = $cell.Value2
    if ($null -ne $current -and $map.ContainsKey($current)) {
        $cell.Value2 = $map[$current]
    }
}
